$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before row 17; this shifts existing rows 17..108 down to 18..109
$ws.Rows.Item(17).Insert()

# Populate the freshly inserted row 17 with the new record
$ws.Cells.Item(17, 1).Value = 11
$ws.Cells.Item(17, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(17, 3).Value = "Bíobío"
$ws.Cells.Item(17, 4).Value = (Get-Date -Year 2022 -Month 5 -Day 26 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(17, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(17, 5).Value = 8
$ws.Cells.Item(17, 6).Value = 100112021
$ws.Cells.Item(17, 7).Value = "Ají"
$ws.Cells.Item(17, 8).Value = "Inferno"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 50
$ws.Cells.Item(17, 11).Value = 18000
$ws.Cells.Item(17, 12).Value = 20000
$ws.Cells.Item(17, 13).Value = 19200
$ws.Cells.Item(17, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(17, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(17, 16).Value = 1600
$ws.Cells.Item(17, 17).Value = 12
$ws.Cells.Item(17, 18).Value = "Hortaliza"
